# Logged Week 15 and simulated Week 16
# Update the "H" (home) row totals on both the OFF and DEF sheets to
# reflect the newly logged/simulated week's cumulative stats.

$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 252
$wsOff.Range("C2").Value = 169
$wsOff.Range("D2").Value = 51
$wsOff.Range("E2").Value = 26
$wsOff.Range("F2").Value = 6

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 232
$wsDef.Range("C2").Value = 154
$wsDef.Range("D2").Value = 59
$wsDef.Range("E2").Value = 25
$wsDef.Range("F2").Value = 4
